$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old columns B:H entirely (clears old values/shared-strings)
$ws.Columns("B:H").Delete()

# Set the new column widths for B..I (A keeps its existing width)
$ws.Columns("B").ColumnWidth = 9.140625
$ws.Columns("C").ColumnWidth = 40.42578125
$ws.Columns("D").ColumnWidth = 20.140625
$ws.Columns("E").ColumnWidth = 12.140625
$ws.Columns("F").ColumnWidth = 29.140625
$ws.Columns("G").ColumnWidth = 11.28515625
$ws.Columns("H").ColumnWidth = 23.42578125
$ws.Columns("I").ColumnWidth = 51.85546875

# Re-populate all the cell values for the new layout
$ws.Range("B1").Value = "IdP"
$ws.Range("C1").Value = "Login Button position"
$ws.Range("D1").Value = "Login button form"
$ws.Range("E1").Value = "Auth form"
$ws.Range("F1").Value = "Cannot detect login button"
$ws.Range("G1").Value = "token/code"
$ws.Range("H1").Value = "Oracle"
$ws.Range("I1").Value = "Additionals"
$ws.Range("A2").Value = "instructables.com"
$ws.Range("C2").Value = "Homepage"
$ws.Range("D2").Value = "a"
$ws.Range("E2").Value = "Popup"
$ws.Range("G2").Value = "vul"
$ws.Range("H2").Value = "Yes"
$ws.Range("A3").Value = "answers.com"
$ws.Range("C3").Value = "homepage"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "Popup"
$ws.Range("G3").Value = "vul"
$ws.Range("H3").Value = "Yes"
$ws.Range("A4").Value = "ehow.com"
$ws.Range("C4").Value = "homepage"
$ws.Range("D4").Value = "span"
$ws.Range("E4").Value = "Popup"
$ws.Range("G4").Value = "sr"
$ws.Range("H4").Value = "Yes"
$ws.Range("A5").Value = "huffingtonpost.com"
$ws.Range("B5").Value = "G+FB"
$ws.Range("C5").Value = "homepage, need to make visible by clicking"
$ws.Range("D5").Value = "a"
$ws.Range("E5").Value = "Popup"
$ws.Range("I5").Value = "import info"
$ws.Range("A6").Value = "pinterest"
$ws.Range("C6").Value = "navigates to signin page"
$ws.Range("D6").Value = "span"
$ws.Range("E6").Value = "navigate"
$ws.Range("I6").Value = "import info"
$ws.Range("A7").Value = "hark.com"
$ws.Range("C7").Value = "navigates to sign up page"
$ws.Range("D7").Value = "a"
$ws.Range("E7").Value = "navigate"
$ws.Range("G7").Value = "code"
$ws.Range("A8").Value = "http://www.squidoo.com/"
$ws.Range("C8").Value = "navigates to signin page"
$ws.Range("D8").Value = "a"
$ws.Range("E8").Value = "Popup"
$ws.Range("G8").Value = "token"
$ws.Range("H8").Value = "Yes"
$ws.Range("I8").Value = "need to agree to user agreement"
$ws.Range("A9").Value = "hulu.com"
$ws.Range("C9").Value = "homepage, need to make visible by clicking"
$ws.Range("E9").Value = "Popup"
$ws.Range("I9").Value = "import info"
$ws.Range("A10").Value = "espn.go.com"
$ws.Range("C10").Value = "homepage"
$ws.Range("D10").Value = "a"
$ws.Range("E10").Value = "popup"
$ws.Range("I10").Value = "import info"
$ws.Range("A11").Value = "ask.com"
$ws.Range("C11").Value = "homepage, need to make visible by clicking"
$ws.Range("D11").Value = "a"
$ws.Range("E11").Value = "popup"
$ws.Range("I11").Value = "import info"
$ws.Range("A12").Value = "imgur.com"
$ws.Range("B12").Value = "G+FB"
$ws.Range("C12").Value = "homepage, need to make visible by clicking"
$ws.Range("I12").Value = "import info"
$ws.Range("A13").Value = "foxnews.com"
$ws.Range("B13").Value = "G+FB"
$ws.Range("C13").Value = "homepage, need to make visible by clicking"
$ws.Range("I13").Value = "import info"
$ws.Range("A14").Value = "flickr.com"
$ws.Range("B14").Value = "G+FB"
$ws.Range("C14").Value = "popup, and then click"
$ws.Range("I14").Value = "import info"
$ws.Range("A15").Value = "cnet.com"
$ws.Range("B15").Value = "G+FB"
$ws.Range("C15").Value = "homepage, need to make visible by clicking"
$ws.Range("I15").Value = "import info"
$ws.Range("A16").Value = "indeed.com"
$ws.Range("A27").Value = "imdb.com"
$ws.Range("B27").Value = "G+FB"
$ws.Range("F27").Value = "Bad, reason: string doesn't have login pattern, just have FB"
$ws.Range("A28").Value = "nbcnews.com"
$ws.Range("F28").Value = "Bad, reason: the login element located is not 'a', clicking its parent node/children would help"
$ws.Range("A33").Value = "netflix.com"
$ws.Range("I33").Value = "Only having a netflix account first and then link it works."
$ws.Range("A34").Value = "zillow.com"
$ws.Range("I34").Value = "import info"
$ws.Range("A35").Value = "reference.com"
$ws.Range("A36").Value = "groupon.com"
$ws.Range("I36").Value = "import info"
$ws.Range("A37").Value = "wikia.com"
$ws.Range("I37").Value = "import info"
$ws.Range("A38").Value = "washingtonpost.com"
$ws.Range("A39").Value = "usatoday.com"
$ws.Range("A40").Value = "vimeo.com"
$ws.Range("I40").Value = "requires linking"
$ws.Range("A41").Value = "dailymail.co.uk"
$ws.Range("I41").Value = "import info"
$ws.Range("A42").Value = "shopathome.com"
$ws.Range("I42").Value = "import info"
$ws.Range("A43").Value = "tripadvisor.com"
$ws.Range("A44").Value = "forbes.com"
$ws.Range("B44").Value = "G+FB"
$ws.Range("A45").Value = "match.com"
$ws.Range("I45").Value = "import info"
$ws.Range("A46").Value = "monster.com"
$ws.Range("I46").Value = "import info"
$ws.Range("A47").Value = "photobucket.com"
$ws.Range("I47").Value = "import info"
$ws.Range("A48").Value = "stackoverflow.com"
$ws.Range("A49").Value = "wsj.com"
$ws.Range("I49").Value = "import info"
$ws.Range("A50").Value = "mapquest.com"
$ws.Range("A51").Value = "swagbucks.com"
$ws.Range("A52").Value = "sears.com"
$ws.Range("B52").Value = "G+FB"
$ws.Range("I52").Value = "import info"
$ws.Range("A53").Value = "buzzfeed.com"
$ws.Range("I53").Value = "logs in automatically, but needs 10 secs or so"
$ws.Range("A54").Value = "yellowpages.com"
$ws.Range("A55").Value = "hootsuite.com"
$ws.Range("A56").Value = "expedia.com"
$ws.Range("I56").Value = "import info"
$ws.Range("A57").Value = "trulia.com"
$ws.Range("A58").Value = "dailymotion.com"
$ws.Range("A59").Value = "careerbuilder.com"
$ws.Range("I59").Value = "import info"
$ws.Range("A60").Value = "pogo.com"
$ws.Range("I60").Value = "need linking first"
$ws.Range("A61").Value = "goodreads.com"
$ws.Range("A62").Value = "foxsports.com"
$ws.Range("I62").Value = "import info"
$ws.Range("A63").Value = "latimes.com"
$ws.Range("I63").Value = "import info"
$ws.Range("A64").Value = "ancestry.com"
$ws.Range("I64").Value = "import info"
$ws.Range("A65").Value = "fiverr.com"
$ws.Range("I65").Value = "import info"
$ws.Range("A66").Value = "slickdeals.net"
$ws.Range("I66").Value = "import info"
$ws.Range("A67").Value = "meetup.com"
$ws.Range("A68").Value = "ign.com"
$ws.Range("I68").Value = "import info"
$ws.Range("A69").Value = "livingsocial.com"
$ws.Range("A70").Value = "coupons.com"
$ws.Range("I70").Value = "need linking first"
$ws.Range("A71").Value = "priceline.com"
$ws.Range("B71").Value = "G+FB"
$ws.Range("I71").Value = "need linking first"
$ws.Range("A72").Value = "reuters.com"
$ws.Range("B72").Value = "G+FB"
$ws.Range("I72").Value = "import info"
$ws.Range("A73").Value = "stumbleupon.com"
$ws.Range("I73").Value = "need linking first"
$ws.Range("A74").Value = "retailmenot.com"
$ws.Range("B74").Value = "G+FB"
$ws.Range("A75").Value = "foodnetwork.com"
$ws.Range("I75").Value = "import info"
$ws.Range("A76").Value = "whitepages.com"
$ws.Range("A77").Value = "wunderground.com"
$ws.Range("I77").Value = "import info"
$ws.Range("A78").Value = "surveymonkey.com"
$ws.Range("I78").Value = "need linking"
$ws.Range("A79").Value = "soundcloud.com"

# Update selection to match the saved view state
$ws.Range("C12").Select()
